$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Buying Opportunity) and Column C (support Zone) updates ---
$ws.Range("B2").Value = "NSE:ADANIENT"
$ws.Range("C2").Value = "NSE:AARTIDRUGS"

$ws.Range("B3").Value = "NSE:APOLLOHOSP"
$ws.Range("C3").Value = "NSE:ADVANIHOTR"

$ws.Range("B4").Value = "NSE:BAJAJFINSV"
$ws.Range("C4").Value = "NSE:AGSTRA"

$ws.Range("B5").Value = "NSE:BASF"
$ws.Range("C5").Value = "NSE:ASTEC"

$ws.Range("B6").Value = "NSE:BPCL"
$ws.Range("C6").Value = "NSE:ASTRAZEN"

$ws.Range("B7").Value = "NSE:BRIGADE"
$ws.Range("C7").Value = "NSE:CDSL"

$ws.Range("B8").Value = "NSE:BSLNIFTY"
$ws.Range("C8").Value = "NSE:CHEMBOND"

$ws.Range("B9").Value = "NSE:COALINDIA"
$ws.Range("C9").Value = "NSE:COMPUSOFT"

$ws.Range("B10").Value = "NSE:FIEMIND"
$ws.Range("C10").Value = "NSE:DIXON"

$ws.Range("B11").Value = "NSE:GRASIM"
$ws.Range("C11").Value = "NSE:GEECEE"

$ws.Range("B12").Value = "NSE:HINDPETRO"
$ws.Range("C12").Value = "NSE:GENUSPAPER"

$ws.Range("B13").Value = "NSE:HONDAPOWER"
$ws.Range("C13").Value = "NSE:GNA"

$ws.Range("B14").Value = "NSE:HPIL"
$ws.Range("C14").Value = "NSE:GRMOVER"

$ws.Range("B15").Value = "NSE:M&M"
$ws.Range("C15").Value = "NSE:GRSE"

$ws.Range("B16").Value = "NSE:MAHKTECH"
$ws.Range("C16").Value = "NSE:HARIOMPIPE"

$ws.Range("B17").Value = "NSE:NITIRAJ"
$ws.Range("C17").Value = "NSE:IRMENERGY"

$ws.Range("B18").Value = "NSE:NKIND"
$ws.Range("C18").Value = "NSE:KERNEX"

$ws.Range("B19").Value = "NSE:PEL"
$ws.Range("C19").Value = "NSE:MCX"

$ws.Range("B20").Value = "NSE:POLYCAB"
$ws.Range("C20").Value = "NSE:PARACABLES"

$ws.Range("B21").Value = "NSE:RELIANCE"
$ws.Range("C21").Value = "NSE:PENIND"

$ws.Range("C22").Value = "NSE:PILANIINVS"
$ws.Range("C23").Value = "NSE:PIXTRANS"
$ws.Range("C24").Value = "NSE:RAJMET"
$ws.Range("C25").Value = "NSE:RUCHIRA"

# --- Row 2 lost its D/E/F entries; row 3 lost its E entry ---
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("E3").ClearContents()

# --- Remove now-unused trailing rows 26-39 (data trimmed to 24 entries) ---
$ws.Range("A26:A39").EntireRow.Delete()
